$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.977669497583861
$ws.Range("J2").Value = 0.977669497583861
$ws.Range("M2").Value = 2.815739333333333
$ws.Range("N2").Value = 8.447217999999999
$ws.Range("O2").Value = 0.07700398964630729
$ws.Range("P2").Value = 0.07700398964630729
$ws.Range("Q2").Value = 21.15970386343666
$ws.Range("R2").Value = 190.43733477093
$ws.Range("S2").Value = 0.07528445186945809
$ws.Range("T2").Value = 0.07528445186945809

$ws.Range("I3").Value = 0.977669497583861
$ws.Range("J3").Value = 0.977669497583861
$ws.Range("O3").Value = 0.1324338085883186
$ws.Range("P3").Value = 0.1324338085883186
$ws.Range("S3").Value = 0.1294764951056586
$ws.Range("T3").Value = 0.1294764951056586

$ws.Range("I4").Value = 0.977669497583861
$ws.Range("J4").Value = 0.977669497583861
$ws.Range("M4").Value = 5.537790999999999
$ws.Range("N4").Value = 16.613373
$ws.Range("O4").Value = 0.1514458372546134
$ws.Range("P4").Value = 0.1514458372546134
$ws.Range("Q4").Value = 41.61536411784499
$ws.Range("R4").Value = 374.538277060605
$ws.Range("S4").Value = 0.1480639756198851
$ws.Range("T4").Value = 0.1480639756198851

$ws.Range("I5").Value = 0.977669497583861
$ws.Range("J5").Value = 0.977669497583861
$ws.Range("M5").Value = 1.188595666666667
$ws.Range("N5").Value = 3.565787
$ws.Range("O5").Value = 0.03250535563648733
$ws.Range("P5").Value = 0.03250535563648733
$ws.Range("Q5").Value = 8.932052772888333
$ws.Range("R5").Value = 80.38847495599499
$ws.Range("S5").Value = 0.03177949471390929
$ws.Range("T5").Value = 0.03177949471390929

$ws.Range("I6").Value = 0.977669497583861
$ws.Range("J6").Value = 0.977669497583861
$ws.Range("M6").Value = 18.85109966666667
$ws.Range("N6").Value = 56.553299
$ws.Range("O6").Value = 0.5155341854158992
$ws.Range("P6").Value = 0.5155341854158992
$ws.Range("Q6").Value = 141.6621495195683
$ws.Range("R6").Value = 1274.959345676115
$ws.Range("S6").Value = 0.5040220480428672
$ws.Range("T6").Value = 0.5040220480428672

$ws.Range("I7").Value = 0.977669497583861
$ws.Range("J7").Value = 0.977669497583861
$ws.Range("M7").Value = 3.330328666666666
$ws.Range("N7").Value = 9.990985999999999
$ws.Range("O7").Value = 0.09107682345837424
$ws.Range("P7").Value = 0.09107682345837424
$ws.Range("Q7").Value = 25.02673721262333
$ws.Range("R7").Value = 225.24063491361
$ws.Range("S7").Value = 0.08904303223208275
$ws.Range("T7").Value = 0.08904303223208275

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.171642
$ws.Range("H8").Value = 0.514926
$ws.Range("I8").Value = 0.02233050241613897
$ws.Range("J8").Value = 0.02233050241613898
$ws.Range("M8").Value = 2.815739333333333
$ws.Range("N8").Value = 8.447217999999999
$ws.Range("O8").Value = 0.07700398964630729
$ws.Range("P8").Value = 0.07700398964630729
$ws.Range("Q8").Value = 0.4832991306519999
$ws.Range("R8").Value = 4.349692175867999
$ws.Range("S8").Value = 0.001719537776849206
$ws.Range("T8").Value = 0.001719537776849206

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.171642
$ws.Range("H9").Value = 0.514926
$ws.Range("I9").Value = 0.02233050241613897
$ws.Range("J9").Value = 0.02233050241613898
$ws.Range("O9").Value = 0.1324338085883186
$ws.Range("P9").Value = 0.1324338085883186
$ws.Range("Q9").Value = 0.831192576562
$ws.Range("R9").Value = 7.480733189057999
$ws.Range("S9").Value = 0.002957313482659934
$ws.Range("T9").Value = 0.002957313482659935

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.171642
$ws.Range("H10").Value = 0.514926
$ws.Range("I10").Value = 0.02233050241613897
$ws.Range("J10").Value = 0.02233050241613898
$ws.Range("M10").Value = 5.537790999999999
$ws.Range("N10").Value = 16.613373
$ws.Range("O10").Value = 0.1514458372546134
$ws.Range("P10").Value = 0.1514458372546134
$ws.Range("Q10").Value = 0.9505175228219999
$ws.Range("R10").Value = 8.554657705398
$ws.Range("S10").Value = 0.003381861634728335
$ws.Range("T10").Value = 0.003381861634728336

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.171642
$ws.Range("H11").Value = 0.514926
$ws.Range("I11").Value = 0.02233050241613897
$ws.Range("J11").Value = 0.02233050241613898
$ws.Range("M11").Value = 1.188595666666667
$ws.Range("N11").Value = 3.565787
$ws.Range("O11").Value = 0.03250535563648733
$ws.Range("P11").Value = 0.03250535563648733
$ws.Range("Q11").Value = 0.204012937418
$ws.Range("R11").Value = 1.836116436762
$ws.Range("S11").Value = 0.0007258609225780368
$ws.Range("T11").Value = 0.0007258609225780369

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.171642
$ws.Range("H12").Value = 0.514926
$ws.Range("I12").Value = 0.02233050241613897
$ws.Range("J12").Value = 0.02233050241613898
$ws.Range("M12").Value = 18.85109966666667
$ws.Range("N12").Value = 56.553299
$ws.Range("O12").Value = 0.5155341854158992
$ws.Range("P12").Value = 0.5155341854158992
$ws.Range("Q12").Value = 3.235640448986
$ws.Range("R12").Value = 29.120764040874
$ws.Range("S12").Value = 0.01151213737303198
$ws.Range("T12").Value = 0.01151213737303198

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.171642
$ws.Range("H13").Value = 0.514926
$ws.Range("I13").Value = 0.02233050241613897
$ws.Range("J13").Value = 0.02233050241613898
$ws.Range("M13").Value = 3.330328666666666
$ws.Range("N13").Value = 9.990985999999999
$ws.Range("O13").Value = 0.09107682345837424
$ws.Range("P13").Value = 0.09107682345837424
$ws.Range("Q13").Value = 0.571624273004
$ws.Range("R13").Value = 5.144618457036
$ws.Range("S13").Value = 0.002033791226291489
$ws.Range("T13").Value = 0.002033791226291489

